# Insert a new weekly data row for Ciboulette (Femacal de La Calera) before the
# existing row 196, shifting the rest of the table down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 196 (existing rows 196:204 shift down to 197:205)
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Cells.Item(196, 1).Value = 3
$ws.Cells.Item(196, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(196, 3).Value = "Coquimbo"
$ws.Cells.Item(196, 4).Value = 44509
$ws.Cells.Item(196, 5).Value = 5
$ws.Cells.Item(196, 6).Value = 100112039
$ws.Cells.Item(196, 7).Value = "Ciboulette"
$ws.Cells.Item(196, 8).Value = "Sin especificar"
$ws.Cells.Item(196, 9).Value = "Primera"
$ws.Cells.Item(196, 10).Value = 160
$ws.Cells.Item(196, 11).Value = 1500
$ws.Cells.Item(196, 12).Value = 1500
$ws.Cells.Item(196, 13).Value = 1500
$ws.Cells.Item(196, 14).Value = "$/docena de atados"
$ws.Cells.Item(196, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(196, 16).Value = 500
$ws.Cells.Item(196, 17).Value = 3
$ws.Cells.Item(196, 18).Value = "Hortaliza"

# Apply the same date style (numFmt "YYYY-MM-DD HH:MM:SS") used on column D elsewhere.
$ws.Cells.Item(196, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
